$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style from H1 (bold, border, centered) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-25
$values = @(
    @(6, 7),
    @(9, 9),
    @(9, 9),
    @(7, 9),
    @(7, 8),
    @(7, 8),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(3, 3),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(3, 3),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
